# "add more to data set"
# Append four new test rows (T0010-T0013) of streaming QoS data to Sheet1,
# matching the formatting already used by neighbouring rows, fix up a couple
# of pre-existing cells whose alignment was out of step with the rest of the
# column, and leave the selection where the author left off (K17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Row 10 - T0010
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = 45994.413888888892

$ws.Range("D10").Value = "Thelakataha Gaatha"
$ws.Range("E10").Value = "U001"
$ws.Range("F10").Value = "Android"
$ws.Range("G10").Value = "9.0.98.1187"
$ws.Range("H10").Value = "Wi-Fi"

$ws.Range("I6").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 10.7

$ws.Range("J6").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J10").Value = 153.595

$ws.Range("K6").Copy()
$ws.Range("K10").PasteSpecial(-4122)
$ws.Range("K10").Value = "Dialog"

$ws.Range("L10").Value = "colombo"
$ws.Range("M10").Value = "Morning"
$ws.Range("N10").Value = 310
$ws.Range("O10").Value = 1.88
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = "play well"

# ---------------------------------------------------------------------------
# Row 11 - T0011
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = 45994.427083333336

$ws.Range("D11").Value = "Bones"
$ws.Range("E11").Value = "U001"
$ws.Range("F11").Value = "Android"
$ws.Range("G11").Value = "9.0.98.1187"

$ws.Range("I6").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = "3G"

$ws.Range("J6").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = 0.03

$ws.Range("J11").Value = 382.795
$ws.Range("J11").HorizontalAlignment = -4152
$ws.Range("J11").VerticalAlignment = -4108
$ws.Range("J11").WrapText = $true

$ws.Range("K6").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("K11").Value = "Dialog"

$ws.Range("L11").Value = "colombo"
$ws.Range("M11").Value = "Morning"
$ws.Range("N11").Value = 165
$ws.Range("O11").Value = 0.98
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = "play well"

# ---------------------------------------------------------------------------
# Row 12 - T0012
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = 45994.46875

$ws.Range("D12").Value = "Ape gam gode"
$ws.Range("E12").Value = "U001"
$ws.Range("F12").Value = "Android"
$ws.Range("G12").Value = "9.0.98.1187"
$ws.Range("H12").Value = "3G"

$ws.Range("J6").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 0.04

$ws.Range("J6").Copy()
$ws.Range("J12").PasteSpecial(-4122)
$ws.Range("J12").Value = 270.464

$ws.Range("K6").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").Value = "Dialog"

$ws.Range("L12").Value = "colombo"
$ws.Range("M12").Value = "Morning"
$ws.Range("N12").Value = 327
$ws.Range("O12").Value = 1.55
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = "play well"

# ---------------------------------------------------------------------------
# Row 13 - T0013
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = 45994.479166666664

$ws.Range("D13").Value = "Demons"
$ws.Range("E13").Value = "U001"
$ws.Range("F13").Value = "Android"
$ws.Range("G13").Value = "9.0.98.1187"
$ws.Range("H13").Value = "4G"

$ws.Range("I6").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = 1.08

$ws.Range("J6").Copy()
$ws.Range("J13").PasteSpecial(-4122)
$ws.Range("J13").Value = 72.428

$ws.Range("K6").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$ws.Range("K13").Value = "Dialog"

$ws.Range("L13").Value = "colombo"
$ws.Range("M13").Value = "Morning"
$ws.Range("N13").Value = 177
$ws.Range("O13").Value = 0.43
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = "play well"

# ---------------------------------------------------------------------------
# Small formatting touch-ups: a couple of cells in column K had fallen out
# of sync with the centered look used everywhere else - nudge them back in
# line while we're in here.
# ---------------------------------------------------------------------------

# K1 (header) should just be horizontally centered, like K2:K4.
$ws.Range("K2").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# K5 should match the centered/wrapped look already used by K6:K9.
$ws.Range("K6").Copy()
$ws.Range("K5").PasteSpecial(-4122)

# The blank trailing rows should be centered/wrapped too, matching K6:K9.
$ws.Range("K6").Copy()
$ws.Range("K14:K25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Leave the cursor where the author left it.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("K17").Select()
